# Update "想去人数" (want-to-go count) figures for several events.
# These events appear both on the "展览" sheet and again on the
# aggregated "全部类型" sheet, so each value must be updated in both
# places to keep the two sheets in sync.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 6;  Old = 496;  New = 497 },
    @{ Row = 13; Old = 6254; New = 6257 },
    @{ Row = 15; Old = 324;  New = 325 },
    @{ Row = 16; Old = 2391; New = 2394 },
    @{ Row = 17; Old = 120;  New = 121 },
    @{ Row = 18; Old = 205;  New = 206 },
    @{ Row = 20; Old = 472;  New = 473 }
)

$rowMap = @{
    6  = 8
    13 = 16
    15 = 19
    16 = 20
    17 = 21
    18 = 22
    20 = 24
}

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

foreach ($u in $updates) {
    $wsExhibition.Range("F$($u.Row)").Value = $u.New

    $allRow = $rowMap[$u.Row]
    $wsAll.Range("F$allRow").Value = $u.New
}
